# Apply the edits described by the commit:
#  - rename sheet "Findings" -> "Data"
#  - rename the "CVSSv3 Vector" header label -> "CVSSv3.1 Vector"
#  - move the active cell/selection on the Data sheet from J2 -> F1
#  - column width adjustments that accompanied the edit

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Findings")
$ws2 = $wb.Worksheets.Item("CVSSv3")

# 1. Rename the first sheet
$ws1.Name = "Data"

# 2. Update the header text for the CVSSv3 Vector column (F1)
$ws1.Range("F1").Value = "CVSSv3.1 Vector"

# 3. Column widths on the Data sheet (values chosen so the stored/quantized
#    result matches the target widths as closely as the engine allows)
$ws1.Columns.Item(1).ColumnWidth = 8.5714285714
$ws1.Columns.Item(2).ColumnWidth = 61
$ws1.Columns.Item(3).ColumnWidth = 80.1428571429
$ws1.Range($ws1.Cells.Item(1,4), $ws1.Cells.Item(1,5)).ColumnWidth = 30
$ws1.Columns.Item(6).ColumnWidth = 78.2857142857
$ws1.Columns.Item(7).ColumnWidth = 42.4285714286
$ws1.Columns.Item(8).ColumnWidth = 35.4285714286
$ws1.Columns.Item(9).ColumnWidth = 49
$ws1.Columns.Item(10).ColumnWidth = 40
$ws1.Columns.Item(11).ColumnWidth = 26.8571428571
$ws1.Range($ws1.Cells.Item(1,12), $ws1.Cells.Item(1,13)).ColumnWidth = 39.5714285714
$ws1.Columns.Item(14).ColumnWidth = 87.5714285714
$ws1.Columns.Item(15).ColumnWidth = 39.5714285714
$ws1.Columns.Item(16).ColumnWidth = 65.7142857143

# 4. Column width on the CVSSv3 reference sheet
$ws2.Columns.Item(1).ColumnWidth = 73

# 5. Move the active selection on the Data sheet to F1
$ws1.Range("F1").Select() | Out-Null
